$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Email" column (header + value) for sending invoices by email
$ws.Range("H1").Value = "Email"
$ws.Range("H2").Value = "adrianrentea01@gmail.com"

# Match the column's autofit width used for the new Email column
$ws.Columns.Item(8).ColumnWidth = 23.67

# Update the active selection to the new cell, like after typing the new value
$ws.Range("H1:H2").Select() | Out-Null
